# Charging decisions - include battery level in charging decision
# Insert two new columns (BatteryLevel, BatteryMinLevel) in front of the
# PowerInside/PowerOutside columns (currently I, J) and add a new test row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns straddling the H/I boundary - this keeps column G
# (Priority) untouched and pushes FromGrid% (H), PowerInside (I),
# PowerOutside (J) and the two current-formula columns (K, L) two slots to
# the right (H->J, I->K, J->L, K->M, L->N), while also correctly extending
# the existing uniform-width column formatting of G:H onto the two new
# columns (so the new I:J end up sharing H's column width, matching the
# original G:H block instead of reverting to the sheet default width).
$ws.Range("H1:I1").EntireColumn.Insert()

# The insert above left the old "FromGrid%" column (now in J) one slot off
# from where it belongs (H) - move it back into H so the new, still-empty
# I/J columns are exactly where BatteryLevel/BatteryMinLevel belong.
$ws.Range("J1:J16").Cut()
$ws.Range("H1").Select()
$ws.Paste()

# Headers for the two new columns
$ws.Range("I1").Value = "BatteryLevel"
$ws.Range("J1").Value = "BatteryMinLevel"

# BatteryLevel / BatteryMinLevel values for the existing 15 test rows.
# All rows get BatteryLevel = 0 and BatteryMinLevel = 25.
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 25
}

# New test case (row 17): battery level above its minimum while outside is
# connected, so charging should be allowed from outside even though
# PowerInside/PowerOutside stay at their defaults.
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1000000
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = "Outside"
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 75
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = 0
$ws.Range("L17").Formula = "=230*6*3*1000"
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 6000

# Restore the sheet's saved selection/active cell as it was left in the
# authored workbook.
$ws.Range("K25").Select()
